# Stubs & Use Casses
# Adds four new "Use Case" rows (45-48) describing the forgot-password /
# password-reset / verify-code flows and the edit-credit-card-address flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 - forgot password
$ws.Range("B46").Value = "connectcollect.com/forgot_password/"
$ws.Range("C46").Value = "User Email"
$ws.Range("D46").Value = "it check email and then send code to user email address"

# Row 48 - password reset
$ws.Range("B48").Value = "connectcollect.com/password_reset/"

# Row 47 - verify password reset code
$ws.Range("C47").Value = "User Email, varification Code"
$ws.Range("B47").Value = "connectcollect.com/verify_password_reset_code/"

# Row 48 - password reset (continued)
$ws.Range("C48").Value = "User Email, varification Code, New Password"

# Row 45 - edit credit card address info
$ws.Range("B45").Value = "connectcollect.com/editt_credit_card_address_info/"
$ws.Range("C45").Value = "Houser number, Street, Town/City, County, Postcode, Country, User_ID, Address_ID"

# Result columns (reuse existing "true/false" shared string)
$ws.Range("D45").Value = "true/false"
$ws.Range("D47").Value = "true/false"
$ws.Range("D48").Value = "true/false"

# Update the view: scroll/zoom/selection to match the new bottom-of-sheet focus
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("D45").Select()
